$wb = $excel.ActiveWorkbook

$sMetaData  = $wb.Worksheets.Item("metaData")
$sMetaData1 = $wb.Worksheets.Item("metaData1")

# Reorder tabs: "metaData" first, "metaData1" second.
$sMetaData1.Move($null, $sMetaData)

# Handles are position-bound, not name-bound, across a Move() call -- re-fetch
# fresh references before doing anything else with either sheet.
$sMetaData  = $wb.Worksheets.Item("metaData")
$sMetaData1 = $wb.Worksheets.Item("metaData1")

# "metaData1" row 3 gets new lookfrom/lookat/image_file/point_cloud_file values.
$sMetaData1.Range("A3").Value = "[0 -10 790]"
$sMetaData1.Range("B3").Value = "[0 -10 0]"
$sMetaData1.Range("D3").Value = "Rh_Narpa_z790.png"
$sMetaData1.Range("E3").Value = "rh_narpa_z790.txt"

# "metaData" keeps its own data, just note its own selection.
$null = $sMetaData.Range("C7").Select()

# "metaData1" becomes the active/selected tab with a new selection.
$null = $sMetaData1.Range("E12").Select()
